# Apply the cell value updates described by the OOXML diff.
# A leading backtick-apostrophe forces Excel to treat the assigned
# text as a literal string (quote-prefix) instead of auto-converting
# numeric-looking values (e.g. "1.00", "0.602") into numbers; the
# Style reset afterwards clears the quote-prefix formatting flag so
# the cell keeps the workbook default style, matching the source file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'42.842.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "`'  +0.75%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "`'2.284.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "`'  -0.76%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "`'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "`'  -0.19%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "`'315.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "`'  -0.53%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "`'104.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "`'  +0.08%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "`'  -1.32%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "`'  +0.12%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "`'0.602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "`'  -1.64%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "`'39.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "`'  -1.43%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "`'0.0900"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "`'  -0.94%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "`'8.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "`'  +0.43%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "`'  +2.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "`'0.999"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "`'  +3.37%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "`'15.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "`'  -1.17%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "`'  -0.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "`'2.287.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "`'  -0.85%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "`'42.774.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "`'  +0.64%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "`'7.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "`'  -1.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "`'13.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "`'  +20.92%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "`'73.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "`'  +0.85%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "`'  +0.22%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "`'261.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "`'  -5.27%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "`'2.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "`'  -3.42%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "`'  +0.38%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "`'10.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "`'  -0.20%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "`'7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "`'  +20.69%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "`'2.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "`'  -0.08%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "`'22.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "`'  -2.47%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "`'37.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "`'  +3.88%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "`'166.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "`'  +0.84%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "`'0.0872"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "`'  +0.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "`'  -3.77%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "`'  -0.49%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "`'  -2.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "`'4.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "`'  -1.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "`'  -5.51%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "`'3.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "`'  +1.80%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "`'2.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "`'  -3.50%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "`'1.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "`'  +4.70%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "`'0.231"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "`'  +1.70%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "`'  -0.02%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "`'1.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "`'  +0.18%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "`'92.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "`'  -2.30%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "`'12.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "`'  +0.97%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "`'Maker"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "`'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "`'1.732.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "`'  +8.62%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "`'Aave"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "`'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "`'113.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "`'  +0.04%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "`'79.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "`'  -3.50%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "`'8.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "`'  -1.80%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "`'5.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "`'  +0.34%  "
$ws.Range("E51").Style = "Normal"
